$wb = $excel.ActiveWorkbook

# --- ListSubCategoryPage (sheet3): drop the imgFilePath/imgPathLink columns,
#     keep only the scrollData column (B) ---
$wsSub = $wb.Worksheets.Item("ListSubCategoryPage")
$wsSub.Range("A1:A2").ClearContents()
$wsSub.Range("C1:C2").ClearContents()

# --- Make ListSubCategoryPage the active/selected tab (chained after
#     ListCategoryPage previously being active) ---
$wsSub.Activate()
